$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.330.27"
$ws.Range("E2").Value = "  +5.86%  "
$ws.Range("D3").Value = "3.540.99"
$ws.Range("E3").Value = "  +5.63%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'190.05"
$ws.Range("E5").Value = "  +10.07%  "
$ws.Range("D6").Value = "'562.31"
$ws.Range("E6").Value = "  +5.89%  "
$ws.Range("D7").Value = "3.536.82"
$ws.Range("E7").Value = "  +5.58%  "
$ws.Range("D8").Value = "'0.617"
$ws.Range("E8").Value = "  +3.53%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "'0.633"
$ws.Range("E10").Value = "  +4.00%  "
$ws.Range("E11").Value = "  +14.62%  "
$ws.Range("D12").Value = "'55.16"
$ws.Range("E12").Value = "  +3.93%  "
$ws.Range("D13").Value = "'0.0000272"
$ws.Range("E13").Value = "  +6.55%  "
$ws.Range("D14").Value = "'9.36"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").Value = "4.109.08"
$ws.Range("E15").Value = "  +5.87%  "
$ws.Range("D16").Value = "3.549.51"
$ws.Range("E16").Value = "  +6.09%  "
$ws.Range("E17").Value = "  +3.71%  "
$ws.Range("D18").Value = "'18.49"
$ws.Range("E18").Value = "  +5.63%  "
$ws.Range("D19").Value = "67.376.14"
$ws.Range("E19").Value = "  +6.03%  "
$ws.Range("E20").Value = "  +7.14%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  +3.41%  "
$ws.Range("D22").Value = "'427.19"
$ws.Range("E22").Value = "  +14.68%  "
$ws.Range("D23").Value = "'4.11"
$ws.Range("E23").Value = "  +9.86%  "
$ws.Range("D24").Value = "'85.38"
$ws.Range("E24").Value = "  +4.63%  "
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").Value = "'11.10"
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("D27").Value = "'2.91"
$ws.Range("E27").Value = "  +7.74%  "
$ws.Range("D28").Value = "'6.16"
$ws.Range("D29").Value = "'12.27"
$ws.Range("E29").Value = "  +8.44%  "
$ws.Range("D30").Value = "'9.02"
$ws.Range("E30").Value = "  +8.94%  "
$ws.Range("D31").Value = "'30.52"
$ws.Range("E31").Value = "  +5.69%  "
$ws.Range("D32").Value = "'632.15"
$ws.Range("E32").Value = "  -2.12%  "
$ws.Range("D33").Value = "'6.67"
$ws.Range("E33").Value = "  +4.17%  "
$ws.Range("D34").Value = "'11.74"
$ws.Range("E34").Value = "  +4.76%  "
$ws.Range("D35").Value = "'0.111"
$ws.Range("E35").Value = "  +4.82%  "
$ws.Range("D36").Value = "'60.12"
$ws.Range("E36").Value = "  +3.56%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "'38.46"
$ws.Range("E37").Value = "  +4.18%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0812"
$ws.Range("E38").Value = "  +11.73%  "
$ws.Range("E39").Value = "  +18.68%  "
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").Value = "'0.388"
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("D42").Value = "'3.42"
$ws.Range("E42").Value = "  +14.09%  "
$ws.Range("D43").Value = "3.132.73"
$ws.Range("E43").Value = "  +6.92%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("D46").Value = "'2.88"
$ws.Range("E46").Value = "  +10.46%  "
$ws.Range("D47").Value = "'3.35"
$ws.Range("E47").Value = "  +9.94%  "
$ws.Range("D48").Value = "'0.0419"
$ws.Range("E48").Value = "  +5.46%  "
$ws.Range("D49").Value = "'2.78"
$ws.Range("E49").Value = "  +4.20%  "
$ws.Range("E50").Value = "  +5.37%  "
$ws.Range("D51").Value = "'141.60"
$ws.Range("E51").Value = "  +3.65%  "
